$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" parameter row (row 5) was dropped from the
# table; deleting the whole row shifts "pie_threshold_range" up into row 5
# and also drops the now-unused "theta_threshold_range" shared string.
$ws.Rows("5").Delete() | Out-Null

# Update the surviving threshold values to the new dataset.
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 0.8

# The shifted-up "pie_threshold_range" row gets new Min/Max values too.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# That row had inherited the special Times New Roman formatting from the
# deleted row below it; re-apply the plain formatting used by the rest of
# the data rows (copied from a normally-formatted cell) so it matches the
# other rows again.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("B5:C5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the saved selection/active cell of the edited workbook.
$ws.Range("G5").Select() | Out-Null
